$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.538445
$ws.Range("H2").Value = 1.615335
$ws.Range("I2").Value = 0.03371608002174246
$ws.Range("J2").Value = 0.03371608002174246
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.310643
$ws.Range("N2").Value = 0.931929
$ws.Range("O2").Value = 0.1259996917369272
$ws.Range("P2").Value = 0.1259996917369272
$ws.Range("Q2").Value = 0.167264170135
$ws.Range("R2").Value = 1.505377531215
$ws.Range("S2").Value = 0.00424821568931712
$ws.Range("T2").Value = 0.00424821568931712

$ws.Range("G3").Value = 0.538445
$ws.Range("H3").Value = 1.615335
$ws.Range("I3").Value = 0.03371608002174246
$ws.Range("J3").Value = 0.03371608002174246
$ws.Range("O3").Value = 0.5516970693375588
$ws.Range("P3").Value = 0.5516970693375588
$ws.Range("Q3").Value = 0.7323760177233333
$ws.Range("R3").Value = 6.59138415951
$ws.Range("S3").Value = 0.01860106253754593
$ws.Range("T3").Value = 0.01860106253754593

$ws.Range("G4").Value = 0.538445
$ws.Range("H4").Value = 1.615335
$ws.Range("I4").Value = 0.03371608002174246
$ws.Range("J4").Value = 0.03371608002174246
$ws.Range("M4").Value = 0.794615
$ws.Range("N4").Value = 2.383845
$ws.Range("O4").Value = 0.3223032389255139
$ws.Range("P4").Value = 0.3223032389255139
$ws.Range("Q4").Value = 0.427856473675
$ws.Range("R4").Value = 3.850708263075
$ws.Range("S4").Value = 0.01086680179487941
$ws.Range("T4").Value = 0.01086680179487941

$ws.Range("I5").Value = 0.7539416098905094
$ws.Range("J5").Value = 0.7539416098905093
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.310643
$ws.Range("N5").Value = 0.931929
$ws.Range("O5").Value = 0.1259996917369272
$ws.Range("P5").Value = 0.1259996917369272
$ws.Range("Q5").Value = 3.740275192942334
$ws.Range("R5").Value = 33.662476736481
$ws.Range("S5").Value = 0.09499641043384681
$ws.Range("T5").Value = 0.0949964104338468

$ws.Range("I6").Value = 0.7539416098905094
$ws.Range("J6").Value = 0.7539416098905093
$ws.Range("O6").Value = 0.5516970693375588
$ws.Range("P6").Value = 0.5516970693375588
$ws.Range("S6").Value = 0.4159473766282351
$ws.Range("T6").Value = 0.415947376628235

$ws.Range("I7").Value = 0.7539416098905094
$ws.Range("J7").Value = 0.7539416098905093
$ws.Range("M7").Value = 0.794615
$ws.Range("N7").Value = 2.383845
$ws.Range("O7").Value = 0.3223032389255139
$ws.Range("P7").Value = 0.3223032389255139
$ws.Range("Q7").Value = 9.567506019578333
$ws.Range("R7").Value = 86.10755417620501
$ws.Range("S7").Value = 0.2429978228284274
$ws.Range("T7").Value = 0.2429978228284274

$ws.Range("G8").Value = 3.391101666666666
$ws.Range("H8").Value = 10.173305
$ws.Range("I8").Value = 0.2123423100877482
$ws.Range("J8").Value = 0.2123423100877481
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.310643
$ws.Range("N8").Value = 0.931929
$ws.Range("O8").Value = 0.1259996917369272
$ws.Range("P8").Value = 0.1259996917369272
$ws.Range("Q8").Value = 1.053421995038333
$ws.Range("R8").Value = 9.480797955344999
$ws.Range("S8").Value = 0.02675506561376328
$ws.Range("T8").Value = 0.02675506561376327

$ws.Range("G9").Value = 3.391101666666666
$ws.Range("H9").Value = 10.173305
$ws.Range("I9").Value = 0.2123423100877482
$ws.Range("J9").Value = 0.2123423100877481
$ws.Range("O9").Value = 0.5516970693375588
$ws.Range("P9").Value = 0.5516970693375588
$ws.Range("Q9").Value = 4.61247023248111
$ws.Range("R9").Value = 41.51223209233
$ws.Range("S9").Value = 0.1171486301717778
$ws.Range("T9").Value = 0.1171486301717778

$ws.Range("G10").Value = 3.391101666666666
$ws.Range("H10").Value = 10.173305
$ws.Range("I10").Value = 0.2123423100877482
$ws.Range("J10").Value = 0.2123423100877481
$ws.Range("M10").Value = 0.794615
$ws.Range("N10").Value = 2.383845
$ws.Range("O10").Value = 0.3223032389255139
$ws.Range("P10").Value = 0.3223032389255139
$ws.Range("Q10").Value = 2.694620250858333
$ws.Range("R10").Value = 24.251582257725
$ws.Range("S10").Value = 0.06843861430220706
$ws.Range("T10").Value = 0.06843861430220705
